$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------
# Add timing data for "Aashish Sort 4" block (rows 132-136).
# Columns P,Q,R,S get numeric trial values; T,U get "Unmeasureable".
# ---------------------------------------------------------------

# Row 132 (Trial 1)
$ws.Range("P132").Value = 3
$ws.Range("Q132").Value = 11
$ws.Range("R132").Value = 349
$ws.Range("S132").Value = 29788
$ws.Range("T132").Value = "Unmeasureable"
$ws.Range("U132").Value = "Unmeasureable"

# Row 133 (Trial 2)
$ws.Range("P133").Value = 1
$ws.Range("Q133").Value = 6
$ws.Range("R133").Value = 286
$ws.Range("S133").Value = 29873
$ws.Range("T133").Value = "Unmeasureable"
$ws.Range("U133").Value = "Unmeasureable"

# Row 134 (Trial 3)
$ws.Range("P134").Value = 2
$ws.Range("Q134").Value = 3
$ws.Range("R134").Value = 185
$ws.Range("S134").Value = 21541
$ws.Range("T134").Value = "Unmeasureable"
$ws.Range("U134").Value = "Unmeasureable"

# Row 135 (Trial 4)
$ws.Range("P135").Value = 1
$ws.Range("Q135").Value = 3
$ws.Range("R135").Value = 193
$ws.Range("S135").Value = 21286
$ws.Range("T135").Value = "Unmeasureable"
$ws.Range("U135").Value = "Unmeasureable"

# Row 136 (Trial 5 / Average row)
$ws.Range("P136").Value = 1
$ws.Range("Q136").Value = 3
$ws.Range("R136").Value = 188
$ws.Range("S136").Value = 21406
$ws.Range("T136").Value = "Unmeasureable"
$ws.Range("U136").Value = "Unmeasureable"

# ---------------------------------------------------------------
# T132 carried the block's top double-border; once it becomes an
# "Unmeasureable" notice (like T127:U131 above it) that border is
# no longer drawn, so clear it.
# ---------------------------------------------------------------
$ws.Range("T132").Borders.Item(8).LineStyle = -4142

# ---------------------------------------------------------------
# Restore the view state saved with the workbook (scrolled down a
# bit further, new active cell at T139).
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("T139").Select()
